$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column E ("Procent") currently stores a 0-1 fraction; convert every
# data row (E2:E101) to a 0-100 percentage scale.
for ($r = 2; $r -le 101; $r++) {
    $cell = $ws.Cells.Item($r, 5)
    $cell.Value2 = [Math]::Round($cell.Value2 * 100, 2)
}
